$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = "'59.825.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.43%  "
$ws.Range("E2").Style = "Normal"

# --- Row 3 ---
$ws.Range("D3").Value = "'2.421.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.77%  "
$ws.Range("E3").Style = "Normal"

# --- Row 4 ---
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"

# --- Row 5 ---
$ws.Range("D5").Value = "'554.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.32%  "
$ws.Range("E5").Style = "Normal"

# --- Row 6 ---
$ws.Range("D6").Value = "'137.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.74%  "
$ws.Range("E6").Style = "Normal"

# --- Row 7 ---
$ws.Range("E7").Value = "'  +0.00%  "
$ws.Range("E7").Style = "Normal"

# --- Row 8 ---
$ws.Range("E8").Value = "'  +1.18%  "
$ws.Range("E8").Style = "Normal"

# --- Row 9 ---
$ws.Range("E9").Value = "'  +4.65%  "
$ws.Range("E9").Style = "Normal"

# --- Row 10 ---
$ws.Range("D10").Value = "'5.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.33%  "
$ws.Range("E10").Style = "Normal"

# --- Row 11 ---
$ws.Range("D11").Value = "'0.360"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.81%  "
$ws.Range("E11").Style = "Normal"

# --- Row 12 ---
$ws.Range("E12").Value = "'  -2.00%  "
$ws.Range("E12").Style = "Normal"

# --- Row 13 ---
$ws.Range("D13").Value = "'24.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.48%  "
$ws.Range("E13").Style = "Normal"

# --- Row 14 ---
$ws.Range("D14").Value = "'2.850.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.95%  "
$ws.Range("E14").Style = "Normal"

# --- Row 15 ---
$ws.Range("D15").Value = "'59.713.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.35%  "
$ws.Range("E15").Style = "Normal"

# --- Row 16 ---
$ws.Range("E16").Value = "'  +4.56%  "
$ws.Range("E16").Style = "Normal"

# --- Row 17 ---
$ws.Range("D17").Value = "'2.411.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.42%  "
$ws.Range("E17").Style = "Normal"

# --- Row 18 ---
$ws.Range("D18").Value = "'11.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +5.74%  "
$ws.Range("E18").Style = "Normal"

# --- Row 19 ---
$ws.Range("E19").Value = "'  +4.30%  "
$ws.Range("E19").Style = "Normal"

# --- Row 20 ---
$ws.Range("D20").Value = "'334.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.27%  "
$ws.Range("E20").Style = "Normal"

# --- Row 21 ---
$ws.Range("D21").Value = "'6.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.45%  "
$ws.Range("E21").Style = "Normal"

# --- Row 22 ---
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.10%  "
$ws.Range("E22").Style = "Normal"

# --- Row 23 ---
$ws.Range("D23").Value = "'64.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.77%  "
$ws.Range("E23").Style = "Normal"

# --- Row 24 ---
$ws.Range("E24").Value = "'  +0.89%  "
$ws.Range("E24").Style = "Normal"

# --- Row 25 ---
$ws.Range("D25").Value = "'8.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.69%  "
$ws.Range("E25").Style = "Normal"

# --- Row 26 ---
$ws.Range("E26").Value = "'  +0.25%  "
$ws.Range("E26").Style = "Normal"

# --- Row 27 ---
$ws.Range("E27").Value = "'  -1.95%  "
$ws.Range("E27").Style = "Normal"

# --- Row 28 ---
$ws.Range("D28").Value = "'0.0₃0791"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +7.20%  "
$ws.Range("E28").Style = "Normal"

# --- Row 29 ---
$ws.Range("E29").Value = "'  +3.00%  "
$ws.Range("E29").Style = "Normal"

# --- Row 30 ---
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'170.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.03%  "
$ws.Range("E30").Style = "Normal"

# --- Row 31 ---
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "'6.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.62%  "
$ws.Range("E31").Style = "Normal"

# --- Row 32 ---
$ws.Range("D32").Value = "'18.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.39%  "
$ws.Range("E32").Style = "Normal"

# --- Row 33 ---
$ws.Range("E33").Value = "'  +0.75%  "
$ws.Range("E33").Style = "Normal"

# --- Row 35 ---
$ws.Range("E35").Value = "'  +4.91%  "
$ws.Range("E35").Style = "Normal"

# --- Row 36 ---
$ws.Range("D36").Value = "'4.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.53%  "
$ws.Range("E36").Style = "Normal"

# --- Row 37 ---
$ws.Range("E37").Value = "'  +0.02%  "
$ws.Range("E37").Style = "Normal"

# --- Row 38 ---
$ws.Range("E38").Value = "'  -1.27%  "
$ws.Range("E38").Style = "Normal"

# --- Row 39 ---
$ws.Range("D39").Value = "'40.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.45%  "
$ws.Range("E39").Style = "Normal"

# --- Row 40 ---
$ws.Range("D40").Value = "'0.422"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +11.70%  "
$ws.Range("E40").Style = "Normal"

# --- Row 41 ---
$ws.Range("D41").Value = "'312.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +7.06%  "
$ws.Range("E41").Style = "Normal"

# --- Row 42 ---
$ws.Range("D42").Value = "'3.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.71%  "
$ws.Range("E42").Style = "Normal"

# --- Row 43 ---
$ws.Range("D43").Value = "'142.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.56%  "
$ws.Range("E43").Style = "Normal"

# --- Row 44 ---
$ws.Range("E44").Value = "'  +2.45%  "
$ws.Range("E44").Style = "Normal"

# --- Row 45 ---
$ws.Range("E45").Value = "'  +3.97%  "
$ws.Range("E45").Style = "Normal"

# --- Row 46 ---
$ws.Range("B46").Value = "Polygon"
$ws.Range("C46").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D46").Value = "'0.412"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +7.06%  "
$ws.Range("E46").Style = "Normal"

# --- Row 47 ---
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'19.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.13%  "
$ws.Range("E47").Style = "Normal"

# --- Row 48 ---
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.571"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.29%  "
$ws.Range("E48").Style = "Normal"

# --- Row 49 ---
$ws.Range("E49").Value = "'  +2.89%  "
$ws.Range("E49").Style = "Normal"

# --- Row 50 ---
$ws.Range("D50").Value = "'11.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.23%  "
$ws.Range("E50").Style = "Normal"

# --- Row 51 ---
$ws.Range("E51").Value = "'  +4.77%  "
$ws.Range("E51").Style = "Normal"
